# Atualização de bases das ligas, do dia: 25-05-2024 às 15:10
#
# The match with id 7802944 (old row 112, Atletico Ottawa x HFX Wanderers,
# not yet played) was removed from the feed. Deleting that row shifts every
# subsequent row up by one, which also happens to clear the now-unused
# "7802944" entry out of the shared-strings table (Excel recompacts it
# automatically), and Excel keeps the worksheet <dimension> in sync too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale fixture row (previously row 112). Everything below
# (old rows 113-115) shifts up to occupy rows 112-114.
$ws.Rows("112").Delete()

# Column A is a plain 0-based sequence number tied to the row's position,
# not to the underlying fixture, so after the shift it must be restored to
# match the row position instead of inheriting the value from the row that
# moved up.
$ws.Range("A112").Value = 110
$ws.Range("A113").Value = 111
$ws.Range("A114").Value = 112

# Refresh the odds that changed for the three still-unplayed fixtures that
# moved from rows 113-115 into rows 112-114.

# Row 112 (Atletico Ottawa x Forge FC)
$ws.Range("M112").Value = 2.375
$ws.Range("O112").Value = 2.625
$ws.Range("Q112").Value = 1.825
$ws.Range("R112").Value = 1.975

# Row 113 (Vancouver FC x Pacific FC CA)
$ws.Range("M113").Value = 2.8
$ws.Range("N113").Value = 3.4
$ws.Range("O113").Value = 2.15
$ws.Range("Q113").Value = 1.85
$ws.Range("R113").Value = 1.95
$ws.Range("T113").Value = 1.825
$ws.Range("U113").Value = 1.975

# Row 114 (Cavalry FC x Valour FC)
$ws.Range("M114").Value = 1.42
$ws.Range("N114").Value = 3.8
$ws.Range("O114").Value = 6.5
$ws.Range("Q114").Value = 1.925
$ws.Range("R114").Value = 1.875
